$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp text
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 16:04"

# Re-sort country rows whose ranking shifted (set country-name text per row)
$ws.Range("A114").Value = "San Marino"
$ws.Range("A115").Value = "Georgia"
$ws.Range("A136").Value = "Sierra Leona"
$ws.Range("A137").Value = "Republica del Chad"
$ws.Range("A138").Value = "Benin"
$ws.Range("A150").Value = "Togo"
$ws.Range("A151").Value = "Gibraltar"
$ws.Range("A152").Value = "Brunei"
$ws.Range("A153").Value = "Guayana Francesa"
$ws.Range("A205").Value = "Seychelles"
$ws.Range("A206").Value = "Montserrat"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 1293907
$ws.Range("C4").Value = 1284
$ws.Range("E4").Value = 999658
$ws.Range("G4").Value = 70
$ws.Range("H4").Value = 76998
$ws.Range("B17").Value = 57306
$ws.Range("C17").Value = 955
$ws.Range("D17").Value = 17041
$ws.Range("E17").Value = 38366
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 1899
$ws.Range("B44").Value = 9943
$ws.Range("C44").Value = 95
$ws.Range("D44").Value = 2453
$ws.Range("E44").Value = 7281
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 209
$ws.Range("D59").Value = 1615
$ws.Range("E59").Value = 3148
$ws.Range("B63").Value = 3778
$ws.Range("C63").Value = 215
$ws.Range("D63").Value = 472
$ws.Range("E63").Value = 3197
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 109
$ws.Range("E66").Value = 2071
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 16
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 114
$ws.Range("E114").Value = 468
$ws.Range("F114").Value = 4
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 41
$ws.Range("B115").Value = 623
$ws.Range("C115").Value = 8
$ws.Range("D115").Value = 288
$ws.Range("E115").Value = 325
$ws.Range("F115").Value = 6
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 10
$ws.Range("B116").Value = 621
$ws.Range("C116").Value = 14
$ws.Range("D116").Value = 202
$ws.Range("E116").Value = 390
$ws.Range("B136").Value = 257
$ws.Range("C136").Value = 26
$ws.Range("D136").Value = 54
$ws.Range("E136").Value = 186
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 17
$ws.Range("B137").Value = 253
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 50
$ws.Range("E137").Value = 176
$ws.Range("H137").Value = 27
$ws.Range("B138").Value = 242
$ws.Range("C138").Value = 102
$ws.Range("D138").Value = 62
$ws.Range("E138").Value = 178
$ws.Range("H138").Value = 2
$ws.Range("B150").Value = 145
$ws.Range("C150").Value = 10
$ws.Range("D150").Value = 85
$ws.Range("E150").Value = 50
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 10
$ws.Range("B151").Value = 144
$ws.Range("D151").Value = 141
$ws.Range("E151").Value = 3
$ws.Range("F151").Value = 0
$ws.Range("H151").Value = 0
$ws.Range("B152").Value = 141
$ws.Range("D152").Value = 132
$ws.Range("E152").Value = 8
$ws.Range("F152").Value = 2
$ws.Range("B153").Value = 138
$ws.Range("D153").Value = 112
$ws.Range("E153").Value = 25
$ws.Range("H153").Value = 1
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
